$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A17").Value = "Neighborhoods.geojson"
$ws.Range("C17").Value = "https://data.baltimorecity.gov/Neighborhoods/Neighborhoods/5cni-ybar"
$ws.Range("B17").Value = "neighborhood defined by Baltimorecitye as of 2010"

$ws.Range("B18").Select()
